$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

# 05-050316TP: row 34, B34 "0d3cbd5cf9a3bf3ff616ce16adc4567b" -> "22745df37b60cda662d05e96d8a86acc"
$ws.Cells.Item(34, 2).Value = "22745df37b60cda662d05e96d8a86acc"

# 03-030003A: row 46, B46 "5705b26efbab4c9a51253a87ddaf53ae" -> "55aee8d5a55feac135ddf23693395da4"
$ws.Cells.Item(46, 2).Value = "55aee8d5a55feac135ddf23693395da4"

# 05-050203TP: row 159, B159 "17e6f09fd8ea8a8972bc475df817080f" -> "d53e48eed45244df3c9f9b4e47f21ea7"
$ws.Cells.Item(159, 2).Value = "d53e48eed45244df3c9f9b4e47f21ea7"

# 05-050203TC: row 169, B169 "6afcb86346c0f16cac73003425cae14d" -> "18bb69df7987c177eaf4582d9167ccee"
$ws.Cells.Item(169, 2).Value = "18bb69df7987c177eaf4582d9167ccee"

# 03-030007A: row 209, B209 "2a104d1a43f2a8d8b185ee2226a15637" -> "a3302a694c2740123584a4ea75f7ea3d"
$ws.Cells.Item(209, 2).Value = "a3302a694c2740123584a4ea75f7ea3d"

# 03-030008TC: row 280, B280 "c503c1eeb22f98e73b7e63b59a5b395b" -> "30b0212a4d791ab012c500e59d3e1bbc"
$ws.Cells.Item(280, 2).Value = "30b0212a4d791ab012c500e59d3e1bbc"

# 05-050201TC: row 281, B281 "91d6cecafdef3ad37838abc58fd1f3c8" -> "dec44b8d200dea67f6a29a4ffd0a9b9e"
$ws.Cells.Item(281, 2).Value = "dec44b8d200dea67f6a29a4ffd0a9b9e"

# 05-0709-070905BTC: row 419, B419 "afba4ee92bb44bede48ddf483ac24705" -> "2ee5add6736bc97726d8045230c25adb"
$ws.Cells.Item(419, 2).Value = "2ee5add6736bc97726d8045230c25adb"

# 03-030015TP: row 421, B421 "a36f4b7630035ba535f49fc7566d6309" -> "66f9ae0dd57d530c7247feebf7db6f43"
$ws.Cells.Item(421, 2).Value = "66f9ae0dd57d530c7247feebf7db6f43"

# 03-030016A: row 473, B473 "3811e2474dd8b07ac0282aa9b7681586" -> "3f6233748c9d480d537076a8e25cd463"
$ws.Cells.Item(473, 2).Value = "3f6233748c9d480d537076a8e25cd463"

# 03-030003TP: row 496, B496 "2dc83515e5a510db584f7c963e055464" -> "b43aa77ee0b5fe894e9404c6f46f3670"
$ws.Cells.Item(496, 2).Value = "b43aa77ee0b5fe894e9404c6f46f3670"

# 05-050317TC: row 514, B514 "1522a941e7773172e4dd4ad354ab0470" -> "f9868f1b583b25bf519efb645fe3fac7"
$ws.Cells.Item(514, 2).Value = "f9868f1b583b25bf519efb645fe3fac7"

# 05-050317TP: row 524, B524 "929b51ea954a9711462847af84dc8432" -> "5800aeb6242332278198f32197ea6a9f"
$ws.Cells.Item(524, 2).Value = "5800aeb6242332278198f32197ea6a9f"

# 05-050317A: row 666, B666 "abf90ea370bd45b91b48fbc900bc506d" -> "a66e3fa37ea8c48de2616e87fbe1968f"
$ws.Cells.Item(666, 2).Value = "a66e3fa37ea8c48de2616e87fbe1968f"

# 05-050206TP: row 680, B680 "dfc9b3ba408aa959d34138ce25d08e59" -> "e9320df071ce0b7c43caf40dfc64571d"
$ws.Cells.Item(680, 2).Value = "e9320df071ce0b7c43caf40dfc64571d"

# 05-050315TC: row 726, B726 "63c9f9c955a1cd66bf998e68d6445a72" -> "f563b6c7c03985296eca81569c13dfed"
$ws.Cells.Item(726, 2).Value = "f563b6c7c03985296eca81569c13dfed"

# 05-050316A: row 729, B729 "52d45121b8d9764e0fdb39e8ce4c0c5e" -> "8008d5734351601ba2b3c5cdad7c76a5"
$ws.Cells.Item(729, 2).Value = "8008d5734351601ba2b3c5cdad7c76a5"

# 05-050207A: row 733, B733 "4c378edcdadf5352ae31165b2ead8eaa" -> "defded7eaff8ac8f153b42f3600ee957"
$ws.Cells.Item(733, 2).Value = "defded7eaff8ac8f153b42f3600ee957"

# 03-030007TC: row 879, B879 "15d943939bfd0cdd4f54081c6b0d1466" -> "72d5422ec0b924600052fea48f361992"
$ws.Cells.Item(879, 2).Value = "72d5422ec0b924600052fea48f361992"

# 03-030007TP: row 892, B892 "192457e23d98cfd2b513d9468704c260" -> "d71eff2eb12cdcc38fb32fc54da54e03"
$ws.Cells.Item(892, 2).Value = "d71eff2eb12cdcc38fb32fc54da54e03"

# 03-030005TC: row 937, B937 "b67452c104a83d6b55dd039d197d8543" -> "1e4907470dac849fea865d5724b45e52"
$ws.Cells.Item(937, 2).Value = "1e4907470dac849fea865d5724b45e52"

# 03-030005TP: row 939, B939 "98a184254ab5092a28ac8710845b3063" -> "f6e3a6425e9082d2faf41ff12a3357a9"
$ws.Cells.Item(939, 2).Value = "f6e3a6425e9082d2faf41ff12a3357a9"

# 03-030016TC: row 951, B951 "42076be6548696c39ca0ade68eaa9607" -> "890ad4d16169d8e4520289d8c831b869"
$ws.Cells.Item(951, 2).Value = "890ad4d16169d8e4520289d8c831b869"

# 03-030016TP: row 955, B955 "ba5014acca4632c127fc7106cedad4fb" -> "ec8951b0c90004edf34c721157014b9d"
$ws.Cells.Item(955, 2).Value = "ec8951b0c90004edf34c721157014b9d"
